$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.399.13'
$ws.Range('E2').Value = '  +5.84%  '

$ws.Range('D3').Value = '3.280.17'
$ws.Range('E3').Value = '  +1.24%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.13%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '408.03'
$ws.Range('E5').Value = '  +3.39%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '111.07'
$ws.Range('E6').Value = '  +3.61%  '

$ws.Range('D7').Value = '3.272.06'
$ws.Range('E7').Value = '  +1.14%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.561'
$ws.Range('E8').Value = '  -2.29%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '1.00'
$ws.Range('E9').Value = '  -0.01%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.610'
$ws.Range('E10').Value = '  -1.38%  '

$ws.Range('E11').Value = '  +10.35%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '38.13'
$ws.Range('E12').Value = '  -2.30%  '

$ws.Range('E13').Value = '  -0.04%  '

$ws.Range('D14').Value = '3.744.65'
$ws.Range('E14').Value = '  -0.09%  '

$ws.Range('E15').Value = '  -1.42%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '18.81'
$ws.Range('E16').Value = '  -1.60%  '

$ws.Range('D17').Value = '3.286.36'
$ws.Range('E17').Value = '  +1.63%  '

$ws.Range('D18').Value = '60.195.06'
$ws.Range('E18').Value = '  +5.85%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.979'
$ws.Range('E19').Value = '  -5.20%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.33'
$ws.Range('E20').Value = '  -4.90%  '

$ws.Range('E21').Value = '  +4.11%  '

$ws.Range('E22').Value = '  -5.22%  '

$ws.Range('B23').Value = 'InternetComputer(DFINITY)'
$ws.Range('C23').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '12.19'
$ws.Range('E23').Value = '  -5.85%  '

$ws.Range('B24').Value = 'BitcoinCash'
$ws.Range('C24').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '292.71'
$ws.Range('E24').Value = '  -1.22%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '72.59'
$ws.Range('E25').Value = '  -1.91%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.04'
$ws.Range('E26').Value = '  -4.10%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '4.47'
$ws.Range('E27').Value = '  +2.59%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '28.25'
$ws.Range('E28').Value = '  +1.48%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.30'
$ws.Range('E29').Value = '  +0.29%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.169'
$ws.Range('E30').Value = '  +0.14%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.39'
$ws.Range('E31').Value = '  -3.90%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.00'
$ws.Range('E32').Value = '  +0.05%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '11.01'
$ws.Range('E33').Value = '  -3.45%  '

$ws.Range('E34').Value = '  -1.79%  '

$ws.Range('B35').Value = 'InjectiveProtocol'
$ws.Range('C35').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '39.26'
$ws.Range('E35').Value = '  +3.70%  '

$ws.Range('B36').Value = 'Toncoin'
$ws.Range('C36').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.37'
$ws.Range('E36').Value = '  +12.24%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '52.19'
$ws.Range('E37').Value = '  +0.77%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0467'
$ws.Range('E38').Value = '  -3.29%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.996'
$ws.Range('E39').Value = '  -0.27%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.02'
$ws.Range('E40').Value = '  +2.35%  '

$ws.Range('B41').Value = 'LidoDAOToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.26'
$ws.Range('E41').Value = '  -7.43%  '

$ws.Range('B42').Value = 'Monero'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '134.62'
$ws.Range('E42').Value = '  +0.20%  '

$ws.Range('B43').Value = 'Stellar'
$ws.Range('C43').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.118'
$ws.Range('E43').Value = '  -2.09%  '

$ws.Range('B44').Value = 'ARBITRUM'
$ws.Range('C44').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.85'
$ws.Range('E44').Value = '  -1.91%  '

$ws.Range('B45').Value = 'TheGraph'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.275'
$ws.Range('E45').Value = '  -2.63%  '

$ws.Range('B46').Value = 'Celestia'
$ws.Range('C46').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '15.97'
$ws.Range('E46').Value = '  -6.16%  '

$ws.Range('B47').Value = 'NEARProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.69'
$ws.Range('E47').Value = '  -6.72%  '

$ws.Range('B48').Value = 'WEMIXToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.18'
$ws.Range('E48').Value = '  +3.87%  '

$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '20.70'
$ws.Range('E49').Value = '  -6.44%  '

$ws.Range('B50').Value = 'Maker'
$ws.Range('C50').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D50').Value = '2.093.52'
$ws.Range('E50').Value = '  -2.94%  '

$ws.Range('B51').Value = 'RocketPoolETH'
$ws.Range('C51').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D51').Value = '3.616.11'
$ws.Range('E51').Value = '  +1.72%  '
